# Apply cell value updates as described by the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 1.11
$ws.Range("P2").Value = 1.32
# Row 3
$ws.Range("F3").Value = 1.74
$ws.Range("H3").Value = 3.9
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1.93
$ws.Range("K3").Value = 6
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.11
$ws.Range("O3").Value = 1.01
$ws.Range("P3").Value = 1.54
$ws.Range("Q3").Value = 1.01
$ws.Range("R3").Value = 1.18
$ws.Range("S3").Value = 3.65
$ws.Range("T3").Value = 1.03
$ws.Range("U3").Value = 1.03
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 1.01
$ws.Range("X3").Value = 980
$ws.Range("Y3").Value = 980
$ws.Range("Z3").Value = 980
$ws.Range("AB3").Value = 980
$ws.Range("AC3").Value = 980
$ws.Range("AD3").Value = 980
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 980
$ws.Range("AG3").Value = 980
$ws.Range("AH3").Value = 980
$ws.Range("AJ3").Value = 980
$ws.Range("AK3").Value = 980
$ws.Range("AN3").Value = 1000
# Row 4
$ws.Range("H4").Value = 6.2
$ws.Range("N4").Value = 4.2
$ws.Range("O4").Value = 1.19
$ws.Range("T4").Value = 1.76
$ws.Range("U4").Value = 2.04
$ws.Range("X4").Value = 25
$ws.Range("Y4").Value = 32
$ws.Range("AA4").Value = 230
$ws.Range("AC4").Value = 12.5
$ws.Range("AD4").Value = 29
$ws.Range("AE4").Value = 110
$ws.Range("AG4").Value = 11
$ws.Range("AH4").Value = 24
$ws.Range("AI4").Value = 85
$ws.Range("AM4").Value = 130
$ws.Range("AO4").Value = 110
# Row 5
$ws.Range("I5").Value = 2.6
$ws.Range("J5").Value = 1.01
$ws.Range("O5").Value = 1.26
$ws.Range("P5").Value = 1.1
$ws.Range("Q5").Value = 1.26
$ws.Range("R5").Value = 1.09
$ws.Range("S5").Value = 1.26
$ws.Range("V5").Value = 1.78
# Row 6
$ws.Range("F6").Value = 1.13
$ws.Range("J6").Value = 1.43
$ws.Range("N6").Value = 1.11
$ws.Range("P6").Value = 1.28
$ws.Range("R6").Value = 1.09
$ws.Range("T6").Value = 1.03
$ws.Range("U6").Value = 1.03
# Row 7
$ws.Range("H7").Value = 1.38
$ws.Range("J7").Value = 1.09
$ws.Range("P7").Value = 1.1
$ws.Range("R7").Value = 1.09
$ws.Range("V7").Value = 1.59
$ws.Range("W7").Value = 1.28
# Row 8
$ws.Range("F8").Value = 1.81
$ws.Range("H8").Value = 1.01
$ws.Range("N8").Value = 1.11
$ws.Range("P8").Value = 1.28
$ws.Range("R8").Value = 1.09
$ws.Range("T8").Value = 1.03
$ws.Range("U8").Value = 1.03
$ws.Range("V8").Value = 1.34
$ws.Range("W8").Value = 1.53
# Row 9
$ws.Range("J9").Value = 1.09
$ws.Range("N9").Value = 1.3
$ws.Range("P9").Value = 1.3
$ws.Range("R9").Value = 1.09
$ws.Range("T9").Value = 1.03
$ws.Range("U9").Value = 1.03
# Row 10
$ws.Range("H10").Value = 1.09
$ws.Range("N10").Value = 1.3
$ws.Range("P10").Value = 1.3
$ws.Range("R10").Value = 1.09
